$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the answer column (D) for rows 4-19 with the question explanations,
# and complete the category (C) for row 8 which was previously blank.

$ws.Range("D4").Value = "26세 여자환자에서 중학교 때 시작되었다면 만성 변비임을 시사하는 소견이며, 신경쓰면 악화되는 소견은 과민성 장증후군의 가능성이 있습니다.  "
# D4 is highlighted in bold blue to emphasize this particular answer
$ws.Range("D4").Font.Bold = $true
$ws.Range("D4").Font.Color = 12611584

$ws.Range("D5").Value = "변비의 원인 중 기능성이라기 보다는 식습관과 관련된 만성 변비의 가능성이 많습니다.  "
$ws.Range("D6").Value = "만성 변비는 일주일에 3회 미만의 배변 횟수가 진단기준이므로 배변 횟수는 진단에 필수적인 질문입니다. "
$ws.Range("D7").Value = "변비시 보이는 변의 형태 중 Bristol type 1 을 확인할 수 있으며, Britol  type 1, 2는 대장통과시간이 Delayed colon transit 의 가능성이 있습니다. "

$ws.Range("C8").Value = "PI"
$ws.Range("D8").Value = "변비를 진단 기준 준 과도한 힘주기, 잔변감등을 확인하는  중요한 병력입니다. "

$ws.Range("D9").Value = "변비와 동반된 소화기증상을 확인하는 질문입니다. "

$ws.Range("D10").Value = "변비와 동반된 소화기증상을 확인하는 질문입니다. "
# D10 switches from justified text to wrapped text (matches D9's wrap style)
$ws.Range("D10").HorizontalAlignment = 1
$ws.Range("D10").WrapText = $true
$ws.Range("D11").Value = "만성변비에서  Defecation sense 가 현저히 감소 되어 있음을 의미하는 질문입니다. "
$ws.Range("D12").Value = "변비와 관련된 대장암, 치핵등의 항문질환, 과민성 장증후군에 동반된 변비 유무를 확인할 수 있는질문입니다. "
$ws.Range("D13").Value = "대장암과 같은 악성종양에 의한 이차성 변비를  Alarm symptom을 의미하는 중요한 질문입니다. "
$ws.Range("D14").Value = "장기간 만성변비의 자가 치료 방법등을 확인하는 질문입니다. "
$ws.Range("D15").Value = "갑상선 질환, 당뇨병,  Scleroderma등에 등에 의한 이차성 변비와의 감별진단에 중요한 질문입니다. "
$ws.Range("D16").Value = "복부 수술력, 고혈압에 의한 고혈압약제, 당뇨에서 동반되는 이차성 변비의 감별진단를 위한 필수적인 질문입니다.  "
$ws.Range("D17").Value = "대장암과 관련된 위험인자 유무를 확인하는 질문입니다. "
$ws.Range("D18").Value = "변비와 관련된 식이 습관과 생활 습관을 확인할 수 있는 질문으로 섬유질 섭취 부족, 작은 활동량이 변비의 원인임을 시사하는 질문 입니다."
$ws.Range("D19").Value = "대장종양등의 가족력은 대장암에 의한 이차성 변비를 확인할 수 있는 질문입니다."

$ws.Range("D22").Value = "물론 검사를 해보아야 겠지만 일차성 변비의 가능성이 있습니다. 그리고 이차성 변비에 대한 검사를 해야 할 것 같습니다. 너무 걱정하지 마세요."
$ws.Range("D23").Value = "직장수지검사를 먼저 시행하고 혈액검사, 항문기능검사, 대장통과시간등을 검사해야 할 것 같습니다. "

$ws.Range("B26").Value = "만성변비의 흔한 원인 5가지는 식이습관과 활동력 감소 같은 생활 습관, 장운동을 감소시키는 항고혈압제, 소염제, 마약성 진통제 등의 약제, 과민성 대장증후군, 당뇨, 파킨슨 병"
$ws.Range("B27").Value = "등이 있습니다. "

$ws.Range("A23").Select()
